$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7923677563667297
$ws.Range("B1").Value = 1.650888800621033
$ws.Range("C1").Value = 5.183144092559814
$ws.Range("D1").Value = 2.082289695739746
$ws.Range("E1").Value = 1.206410884857178
